# 2040_DK.xlsx — "Updated Results with corrected code"
#
# Changes:
#   - D3 (Hydrogen / Non-metallic minerals) no longer has a value.
#   - Row 7's label changes from "Other" to "Biogas" and its value is
#     corrected to 416.0443856700242.
#   - A new row 8 is inserted with label "Other" and value 1130.709280420575,
#     formatted the same way as the other row labels in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hydrogen / Non-metallic minerals no longer has a value.
$ws.Range("D3").ClearContents()

# "Other" was renamed to "Biogas" and its corrected value filled in.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 416.0443856700242

# A new "Other" row is added below, reusing row 7's label formatting
# (bold, bordered, centered) for the new row 8 label cell.
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 1130.709280420575
